$wb = $excel.ActiveWorkbook

# Fix property_category values that were incorrectly set to "land" in
# the 建物 (building) and 汽車 (car) sheets.

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I16").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H6").Value = "car"
